$d = $word.ActiveDocument

# 1) Remove the stray _GoBack bookmark from the empty paragraph right after the title.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Clean up "Items/power ups" so it is a single run without the spell-check
#    proofErr markers that used to straddle "Items/power " and "ups". A plain
#    Find/Replace merges the two runs' text but leaves a stray <w:proofErr/>
#    behind, so instead build a fresh paragraph (inheriting the list
#    formatting) right after the old one, then delete the old, broken one.
$oldItemsPara = $d.Paragraphs.Item(3)
$oldItemsPara.Range.InsertParagraphAfter()
$newItemsPara = $d.Paragraphs.Item(4)
$newItemsPara.Range.Text = "Items/power ups"
$oldItemsRange = $d.Range($oldItemsPara.Range.Start, $newItemsPara.Range.Start)
$oldItemsRange.Delete()

# 3) Append four new bullet paragraphs after the last existing one, reusing its
#    list formatting (pStyle "Listenabsatz" + numPr ilvl 0 / numId 1), which
#    Word automatically carries over to paragraphs created via InsertParagraphAfter.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.Text = "Cutscene Einstieg, Roboter spawnen, T" + [char]0x00FC + "ren " + [char]0x00F6 + "ffnen sich. Spieler bekommt gesagt: " + [char]0x201E + "Lass sie nicht entkommen!" + [char]0x201C

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "Platforms move"

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "Spawn points move"

$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last

# Give the final paragraph some placeholder text, then delete just the text
# (not the paragraph mark) so no leftover empty run remains, mirroring the
# clean "bookmark only" paragraph from the diff.
$p4.Range.Text = "TEMP"
$textOnly = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$textOnly.Delete()

$finalPara = $d.Paragraphs.Last
$d.Bookmarks.Add("_GoBack", $finalPara.Range)
